# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$texto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.48 = 39052.13 pesos`n✅ 39052.13 pesos = 9.41 = 959.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $texto

# --- tasas: update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 105.5
$wsTasas.Range("O10").Value = 4120
$wsTasas.Range("N12").Value = 4152
$wsTasas.Range("O12").Value = 102
